$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new log entry row 35
$ws.Range("A35").NumberFormat = $ws.Range("A34").NumberFormat
$ws.Range("A35").Value = 45670
$ws.Range("B35").Value = "finshed mansion area"
$ws.Range("C35").Value = 4

# Update selection to match the post-edit state (active cell D36)
$ws.Range("D36").Select()
